# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" values (column G) computed to replace the old Strike# values.
$kValues = @{
    2  = 6
    3  = 10
    4  = 3
    5  = 6
    6  = 4
    7  = 5
    8  = 3
    9  = 4
    10 = 9
    11 = 9
    12 = 6
    13 = 5
    14 = 9
    15 = 5
    16 = 6
    17 = 5
    18 = 8
    19 = 7
    20 = 6
    21 = 6
    22 = 7
    23 = 8
    24 = 9
    25 = 10
    26 = 8
    27 = 10
    28 = 6
    29 = 4
    30 = 5
    31 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
